# Scheduled data refresh: update the computed price/profit columns
# (H:N - currentAveragePrice*, LevePrice*, LeveProfit*) on each leve
# sheet with freshly pulled market-board figures. A few rows' HQ
# profit (or, in one case, NQ cost) cells are blank because that
# leve item has no HQ-crafted variant this cycle, so those cells are
# cleared entirely rather than written as 0.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 622722
$ws.Range("I132").Value = 2243.836
$ws.Range("J132").Value = 4087058.2
$ws.Range("K132").Value = 6731.508
$ws.Range("L132").Value = 12261174.6
$ws.Range("M132").Value = -4201.508
$ws.Range("N132").Value = -12266234.6
$ws.Range("H138").Value = 2566700.5
$ws.Range("I138").Value = 2075.2068
$ws.Range("J138").Value = 4084540
$ws.Range("K138").Value = 6225.6204
$ws.Range("L138").Value = 12253620
$ws.Range("M138").Value = -1085.6204
$ws.Range("N138").Value = -12263900

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3682.82
$ws.Range("I32").Value = 3682.82
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3682.82
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3395.82
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 100202250
$ws.Range("I61").Value = 166836180
$ws.Range("J61").Value = 251375
$ws.Range("K61").Value = 166836180
$ws.Range("L61").Value = 251375
$ws.Range("M61").Value = -166835968
$ws.Range("N61").Value = -251799
$ws.Range("H74").Value = 5325252
$ws.Range("I74").Value = 7115154.5
$ws.Range("J74").Value = 67413.375
$ws.Range("K74").Value = 7115154.5
$ws.Range("L74").Value = 67413.375
$ws.Range("M74").Value = -7114280.5
$ws.Range("N74").Value = -69161.375
$ws.Range("H77").Value = 5325252
$ws.Range("I77").Value = 7115154.5
$ws.Range("J77").Value = 67413.375
$ws.Range("K77").Value = 35575772.5
$ws.Range("L77").Value = 337066.875
$ws.Range("M77").Value = -35571404.5
$ws.Range("N77").Value = -345802.875
$ws.Range("H110").Value = 770593.4
$ws.Range("I110").Value = 2000458
$ws.Range("J110").Value = 1928
$ws.Range("K110").Value = 2000458
$ws.Range("L110").Value = 1928
$ws.Range("M110").Value = -1998413
$ws.Range("N110").Value = -6018
$ws.Range("H132").Value = 7073208.5
$ws.Range("I132").Value = 8081922.5
$ws.Range("K132").Value = 24245767.5
$ws.Range("M132").Value = -24243237.5
$ws.Range("H136").Value = 100202250
$ws.Range("I136").Value = 166836180
$ws.Range("J136").Value = 251375
$ws.Range("K136").Value = 500508540
$ws.Range("L136").Value = 754125
$ws.Range("M136").Value = -500505990
$ws.Range("N136").Value = -759225

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1924.2307
$ws.Range("I20").Value = 1798.2858
$ws.Range("J20").Value = 2071.1667
$ws.Range("K20").Value = 1798.2858
$ws.Range("L20").Value = 2071.1667
$ws.Range("M20").Value = -1551.2858
$ws.Range("N20").Value = -2565.1667
$ws.Range("H107").Value = 4493.4814
$ws.Range("I107").Value = 3806.0588
$ws.Range("K107").Value = 3806.0588
$ws.Range("M107").Value = -1886.0588

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 14086023
$ws.Range("I58").Value = 27028330
$ws.Range("J58").Value = 1747.5
$ws.Range("K58").Value = 27028330
$ws.Range("L58").Value = 1747.5
$ws.Range("M58").Value = -27028127
$ws.Range("N58").Value = -2153.5
$ws.Range("H136").Value = 14086023
$ws.Range("I136").Value = 27028330
$ws.Range("J136").Value = 1747.5
$ws.Range("K136").Value = 81084990
$ws.Range("L136").Value = 5242.5
$ws.Range("M136").Value = -81082440
$ws.Range("N136").Value = -10342.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3450659.5
$ws.Range("I4").Value = 74
$ws.Range("J4").Value = 3848804
$ws.Range("K4").Value = 222
$ws.Range("L4").Value = 11546412
$ws.Range("M4").Value = -110
$ws.Range("N4").Value = -11546636
$ws.Range("H107").Value = 884.2308
$ws.Range("I107").Value = 788.2778
$ws.Range("J107").Value = 1100.125
$ws.Range("K107").Value = 2364.8334
$ws.Range("L107").Value = 3300.375
$ws.Range("M107").Value = -444.8334
$ws.Range("N107").Value = -7140.375
$ws.Range("H131").Value = 820.07574
$ws.Range("J131").Value = 935.7925
$ws.Range("L131").Value = 2807.3775
$ws.Range("N131").Value = -12887.3775
$ws.Range("H132").Value = 2120.818
$ws.Range("I132").Value = 1268.125
$ws.Range("J132").Value = 2923.353
$ws.Range("K132").Value = 11413.125
$ws.Range("L132").Value = 26310.177
$ws.Range("M132").Value = -8883.125
$ws.Range("N132").Value = -31370.177

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000004
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3900
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3900
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3900
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -4124
$ws.Range("H40").Value = 2716
$ws.Range("I40").Value = 2501
$ws.Range("J40").Value = 2802
$ws.Range("K40").Value = 2501
$ws.Range("L40").Value = 2802
$ws.Range("M40").Value = -2365
$ws.Range("N40").Value = -3074
$ws.Range("H45").Value = 15000
$ws.Range("J45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -15814

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 53502.5
$ws.Range("I2").Value = 7002
$ws.Range("J2").Value = 100003
$ws.Range("K2").Value = 7002
$ws.Range("L2").Value = 100003
$ws.Range("M2").Value = -6890
$ws.Range("N2").Value = -100227
$ws.Range("H41").Value = 8320.286
$ws.Range("I41").Value = 8121
$ws.Range("J41").Value = 8400
$ws.Range("K41").Value = 8121
$ws.Range("L41").Value = 8400
$ws.Range("M41").Value = -7731
$ws.Range("N41").Value = -9180
$ws.Range("H62").Value = 4000.5
$ws.Range("I62").Value = 4000.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4000.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3376.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4000.5
$ws.Range("I65").Value = 4000.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20002.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16882.5
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 2222.3333
$ws.Range("I81").Value = 1600.2
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3200.4
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2139.4
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2222.3333
$ws.Range("I84").Value = 1600.2
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 16002
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -10698
$ws.Range("N84").Value = -40608
$ws.Range("H122").Value = 1499.2188
$ws.Range("I122").Value = 1144.2858
$ws.Range("J122").Value = 2176.818
$ws.Range("K122").Value = 3432.8574
$ws.Range("L122").Value = 6530.454000000001
$ws.Range("M122").Value = -982.8574000000003
$ws.Range("N122").Value = -11430.454
$ws.Range("H126").Value = 1364.3
$ws.Range("I126").Value = 659
$ws.Range("J126").Value = 3010
$ws.Range("K126").Value = 1977
$ws.Range("L126").Value = 9030
$ws.Range("M126").Value = 493
$ws.Range("N126").Value = -13970
$ws.Range("H132").Value = 37894.527
$ws.Range("I132").Value = 33319.387
$ws.Range("J132").Value = 43804.082
$ws.Range("K132").Value = 99958.16100000001
$ws.Range("L132").Value = 131412.246
$ws.Range("M132").Value = -97428.16100000001
$ws.Range("N132").Value = -136472.246
$ws.Range("H133").Value = 28738.334
$ws.Range("J133").Value = 28738.334
$ws.Range("L133").Value = 28738.334
$ws.Range("N133").Value = -38858.334
